$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.723.59"
$ws.Range("E2").Value = "  +4.08%  "
$ws.Range("D3").Value = "3.075.47"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'552.40"
$ws.Range("E5").Value = "  +4.48%  "
$ws.Range("D6").Value = "'138.35"
$ws.Range("E6").Value = "  +5.32%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.066.87"
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").Value = "'6.21"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "'0.456"
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("E13").Value = "  +4.43%  "
$ws.Range("D14").Value = "'34.93"
$ws.Range("E14").Value = "  +4.61%  "
$ws.Range("D15").Value = "3.569.25"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "63.571.01"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("D17").Value = "3.076.29"
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("E19").Value = "  +3.82%  "
$ws.Range("D20").Value = "'483.01"
$ws.Range("E20").Value = "  +5.55%  "
$ws.Range("D21").Value = "'13.52"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("D22").Value = "'0.685"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").Value = "'7.21"
$ws.Range("E23").Value = "  +5.53%  "
$ws.Range("D24").Value = "'81.65"
$ws.Range("E24").Value = "  +4.93%  "
$ws.Range("E25").Value = "  +6.25%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E27").Value = "  +4.21%  "
$ws.Range("D28").Value = "'8.01"
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("E29").Value = "  +8.85%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").Value = "'26.04"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("D33").Value = "'5.80"
$ws.Range("E33").Value = "  +7.86%  "
$ws.Range("D34").Value = "'2.41"
$ws.Range("E34").Value = "  +8.07%  "
$ws.Range("D35").Value = "'55.61"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("D37").Value = "'469.09"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("D38").Value = "3.178.95"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "'0.0818"
$ws.Range("E39").Value = "  +5.33%  "
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("D41").Value = "'0.121"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").Value = "'8.23"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'28.24"
$ws.Range("E43").Value = "  +11.31%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.54"
$ws.Range("E44").Value = "  +6.48%  "
$ws.Range("E45").Value = "  +3.47%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "'2.03"
$ws.Range("E47").Value = "  +6.14%  "
$ws.Range("E48").Value = "  +2.30%  "
$ws.Range("D49").Value = "0.0₃0514"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "'116.02"
$ws.Range("E50").Value = "  -3.68%  "
$ws.Range("D51").Value = "'2.06"
$ws.Range("E51").Value = "  +5.61%  "
